# PWR_Board_TestReportTemplate.xlsx edit
# - Insert a new "MEASURED" column into the Quantities sheet (between NOMINAL and SCALE),
#   shifting the existing SCALE/OFFSET columns one column to the right.
# - Switch the active/selected sheet from Report to Quantities, and update the
#   selected cell on Quantities to G20.

$wb = $excel.ActiveWorkbook
$wsQuant = $wb.Worksheets.Item("Quantities")

# Insert a new column before the existing SCALE column (F) on the Quantities sheet.
# This shifts SCALE: F -> G and OFFSET: G -> H, leaving column F blank.
$wsQuant.Range("F1").EntireColumn.Insert()

# Give the new column its header label.
$wsQuant.Range("F1").Value = "MEASURED"

# Update the selection on Quantities, then make Quantities the active sheet/tab
# (this also clears the tabSelected flag on the Report sheet).
$wsQuant.Range("G20").Select()
$wsQuant.Activate()
